$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused columns L..AE (rows 1-3) entirely, including their
# formatting, so the cells disappear from the sheet (dimension shrinks to
# A1:K3).
$ws.Range("L1:AE3").ClearFormats()
$ws.Range("L1:AE3").Value = $null

# Row 2 (B2:K2) -- new distribution values
$row2 = @(
    -0.3999999996542806,
    -0.2999999996926938,
    -0.199999999731107,
    -0.09999999976952034,
    (1.920663628141028 / 10000000000),
    0.1000000001536532,
    0.2000000001152399,
    0.3000000000768266,
    0.4000000000384133,
    0.5
)
for ($i = 0; $i -lt $row2.Count; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
}

# Row 3 (B3:K3) -- new count values
$row3 = @(
    1856,
    1664,
    5596,
    27432,
    694176,
    2338024,
    5470,
    468,
    54,
    22
)
for ($i = 0; $i -lt $row3.Count; $i++) {
    $ws.Cells.Item(3, $i + 2).Value = $row3[$i]
}
